# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E16:E21) previously listed periods in
# ascending order (1909, 1910, 1911, 1912, 2001, 2002). The old periods
# are removed and replaced with the new/current periods, now listed in
# descending (most recent first) order: 2002, 2001, 1912, 1911, 1910, 1909.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("2002", "2001", "1912", "1911", "1910", "1909")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
